$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$s2D = $ws.Range("D2").Style
$s2E = $ws.Range("E2").Style
$ws.Range("D2").Value = "'42.867.64"
$ws.Range("E2").Value = "'  -0.08%  "
$ws.Range("D2").Style = $s2D
$ws.Range("E2").Style = $s2E

# Row 3
$s3D = $ws.Range("D3").Style
$s3E = $ws.Range("E3").Style
$ws.Range("D3").Value = "'2.360.94"
$ws.Range("E3").Value = "'  +1.89%  "
$ws.Range("D3").Style = $s3D
$ws.Range("E3").Style = $s3E

# Row 4
$s4D = $ws.Range("D4").Style
$s4E = $ws.Range("E4").Style
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("D4").Style = $s4D
$ws.Range("E4").Style = $s4E

# Row 5
$s5D = $ws.Range("D5").Style
$s5E = $ws.Range("E5").Style
$ws.Range("D5").Value = "'301.57"
$ws.Range("E5").Value = "'  -0.30%  "
$ws.Range("D5").Style = $s5D
$ws.Range("E5").Style = $s5E

# Row 6
$s6D = $ws.Range("D6").Style
$s6E = $ws.Range("E6").Style
$ws.Range("D6").Value = "'95.20"
$ws.Range("E6").Value = "'  -0.72%  "
$ws.Range("D6").Style = $s6D
$ws.Range("E6").Style = $s6E

# Row 7
$s7D = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.504"
$ws.Range("D7").Style = $s7D

# Row 8
$s8E = $ws.Range("E8").Style
$ws.Range("E8").Value = "'  -0.05%  "
$ws.Range("E8").Style = $s8E

# Row 9
$s9D = $ws.Range("D9").Style
$s9E = $ws.Range("E9").Style
$ws.Range("D9").Value = "'0.485"
$ws.Range("E9").Value = "'  -1.86%  "
$ws.Range("D9").Style = $s9D
$ws.Range("E9").Style = $s9E

# Row 10
$s10D = $ws.Range("D10").Style
$s10E = $ws.Range("E10").Style
$ws.Range("D10").Value = "'33.81"
$ws.Range("E10").Value = "'  -1.80%  "
$ws.Range("D10").Style = $s10D
$ws.Range("E10").Style = $s10E

# Row 11
$s11D = $ws.Range("D11").Style
$s11E = $ws.Range("E11").Style
$ws.Range("D11").Value = "'0.0785"
$ws.Range("E11").Value = "'  +0.07%  "
$ws.Range("D11").Style = $s11D
$ws.Range("E11").Style = $s11E

# Row 12
$s12E = $ws.Range("E12").Style
$ws.Range("E12").Value = "'  +2.69%  "
$ws.Range("E12").Style = $s12E

# Row 13
$s13D = $ws.Range("D13").Style
$s13E = $ws.Range("E13").Style
$ws.Range("D13").Value = "'18.21"
$ws.Range("E13").Value = "'  -4.28%  "
$ws.Range("D13").Style = $s13D
$ws.Range("E13").Style = $s13E

# Row 14
$s14B = $ws.Range("B14").Style
$s14C = $ws.Range("C14").Style
$s14D = $ws.Range("D14").Style
$s14E = $ws.Range("E14").Style
$ws.Range("B14").Value = "'Polkadot"
$ws.Range("C14").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.71"
$ws.Range("E14").Value = "'  -0.52%  "
$ws.Range("B14").Style = $s14B
$ws.Range("C14").Style = $s14C
$ws.Range("D14").Style = $s14D
$ws.Range("E14").Style = $s14E

# Row 15
$s15B = $ws.Range("B15").Style
$s15C = $ws.Range("C15").Style
$s15D = $ws.Range("D15").Style
$s15E = $ws.Range("E15").Style
$ws.Range("B15").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "'2.734.07"
$ws.Range("E15").Value = "'  +2.09%  "
$ws.Range("B15").Style = $s15B
$ws.Range("C15").Style = $s15C
$ws.Range("D15").Style = $s15D
$ws.Range("E15").Style = $s15E

# Row 16
$s16D = $ws.Range("D16").Style
$s16E = $ws.Range("E16").Style
$ws.Range("D16").Value = "'2.353.15"
$ws.Range("E16").Value = "'  +1.74%  "
$ws.Range("D16").Style = $s16D
$ws.Range("E16").Style = $s16E

# Row 17
$s17D = $ws.Range("D17").Style
$s17E = $ws.Range("E17").Style
$ws.Range("D17").Value = "'0.795"
$ws.Range("E17").Value = "'  +0.60%  "
$ws.Range("D17").Style = $s17D
$ws.Range("E17").Style = $s17E

# Row 18
$s18D = $ws.Range("D18").Style
$s18E = $ws.Range("E18").Style
$ws.Range("D18").Value = "'42.808.71"
$ws.Range("E18").Value = "'  -0.03%  "
$ws.Range("D18").Style = $s18D
$ws.Range("E18").Style = $s18E

# Row 19
$s19D = $ws.Range("D19").Style
$s19E = $ws.Range("E19").Style
$ws.Range("D19").Value = "'12.00"
$ws.Range("E19").Value = "'  -1.60%  "
$ws.Range("D19").Style = $s19D
$ws.Range("E19").Style = $s19E

# Row 20
$s20E = $ws.Range("E20").Style
$ws.Range("E20").Value = "'  +1.76%  "
$ws.Range("E20").Style = $s20E

# Row 21
$s21D = $ws.Range("D21").Style
$s21E = $ws.Range("E21").Style
$ws.Range("D21").Value = "'0.0₃0884"
$ws.Range("E21").Value = "'  -1.02%  "
$ws.Range("D21").Style = $s21D
$ws.Range("E21").Style = $s21E

# Row 22
$s22D = $ws.Range("D22").Style
$s22E = $ws.Range("E22").Style
$ws.Range("D22").Value = "'67.83"
$ws.Range("E22").Value = "'  +0.00%  "
$ws.Range("D22").Style = $s22D
$ws.Range("E22").Style = $s22E

# Row 23
$s23D = $ws.Range("D23").Style
$s23E = $ws.Range("E23").Style
$ws.Range("D23").Value = "'234.92"
$ws.Range("E23").Value = "'  -0.31%  "
$ws.Range("D23").Style = $s23D
$ws.Range("E23").Style = $s23E

# Row 24
$s24E = $ws.Range("E24").Style
$ws.Range("E24").Value = "'  -2.31%  "
$ws.Range("E24").Style = $s24E

# Row 25
$s25E = $ws.Range("E25").Style
$ws.Range("E25").Value = "'  -0.04%  "
$ws.Range("E25").Style = $s25E

# Row 26
$s26E = $ws.Range("E26").Style
$ws.Range("E26").Value = "'  -0.11%  "
$ws.Range("E26").Style = $s26E

# Row 27
$s27D = $ws.Range("D27").Style
$s27E = $ws.Range("E27").Style
$ws.Range("D27").Value = "'24.61"
$ws.Range("E27").Value = "'  +0.81%  "
$ws.Range("D27").Style = $s27D
$ws.Range("E27").Style = $s27E

# Row 28
$s28E = $ws.Range("E28").Style
$ws.Range("E28").Value = "'  +0.26%  "
$ws.Range("E28").Style = $s28E

# Row 29
$s29D = $ws.Range("D29").Style
$s29E = $ws.Range("E29").Style
$ws.Range("D29").Value = "'9.23"
$ws.Range("E29").Value = "'  +1.00%  "
$ws.Range("D29").Style = $s29D
$ws.Range("E29").Style = $s29E

# Row 30
$s30D = $ws.Range("D30").Style
$s30E = $ws.Range("E30").Style
$ws.Range("D30").Value = "'31.36"
$ws.Range("E30").Value = "'  -2.86%  "
$ws.Range("D30").Style = $s30D
$ws.Range("E30").Style = $s30E

# Row 31
$s31D = $ws.Range("D31").Style
$s31E = $ws.Range("E31").Style
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "'  -0.04%  "
$ws.Range("D31").Style = $s31D
$ws.Range("E31").Style = $s31E

# Row 32
$s32D = $ws.Range("D32").Style
$s32E = $ws.Range("E32").Style
$ws.Range("D32").Value = "'5.03"
$ws.Range("E32").Value = "'  +0.35%  "
$ws.Range("D32").Style = $s32D
$ws.Range("E32").Style = $s32E

# Row 33
$s33D = $ws.Range("D33").Style
$s33E = $ws.Range("E33").Style
$ws.Range("D33").Value = "'0.0725"
$ws.Range("E33").Value = "'  +3.85%  "
$ws.Range("D33").Style = $s33D
$ws.Range("E33").Style = $s33E

# Row 34
$s34E = $ws.Range("E34").Style
$ws.Range("E34").Value = "'  -3.57%  "
$ws.Range("E34").Style = $s34E

# Row 35
$s35B = $ws.Range("B35").Style
$s35C = $ws.Range("C35").Style
$s35D = $ws.Range("D35").Style
$s35E = $ws.Range("E35").Style
$ws.Range("B35").Value = "'Kaspa"
$ws.Range("C35").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.104"
$ws.Range("E35").Value = "'  +4.33%  "
$ws.Range("B35").Style = $s35B
$ws.Range("C35").Style = $s35C
$ws.Range("D35").Style = $s35D
$ws.Range("E35").Style = $s35E

# Row 36
$s36B = $ws.Range("B36").Style
$s36C = $ws.Range("C36").Style
$s36D = $ws.Range("D36").Style
$s36E = $ws.Range("E36").Style
$ws.Range("B36").Value = "'ARBITRUM"
$ws.Range("C36").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.85"
$ws.Range("E36").Value = "'  +3.28%  "
$ws.Range("B36").Style = $s36B
$ws.Range("C36").Style = $s36C
$ws.Range("D36").Style = $s36D
$ws.Range("E36").Style = $s36E

# Row 37
$s37D = $ws.Range("D37").Style
$s37E = $ws.Range("E37").Style
$ws.Range("D37").Value = "'4.35"
$ws.Range("E37").Value = "'  -2.87%  "
$ws.Range("D37").Style = $s37D
$ws.Range("E37").Style = $s37E

# Row 38
$s38E = $ws.Range("E38").Style
$ws.Range("E38").Value = "'  -1.10%  "
$ws.Range("E38").Style = $s38E

# Row 39
$s39E = $ws.Range("E39").Style
$ws.Range("E39").Value = "'  +1.45%  "
$ws.Range("E39").Style = $s39E

# Row 40
$s40D = $ws.Range("D40").Style
$s40E = $ws.Range("E40").Style
$ws.Range("D40").Value = "'121.22"
$ws.Range("E40").Value = "'  -27.08%  "
$ws.Range("D40").Style = $s40D
$ws.Range("E40").Style = $s40E

# Row 41
$s41E = $ws.Range("E41").Style
$ws.Range("E41").Value = "'  -0.74%  "
$ws.Range("E41").Style = $s41E

# Row 42
$s42D = $ws.Range("D42").Style
$s42E = $ws.Range("E42").Style
$ws.Range("D42").Value = "'21.42"
$ws.Range("E42").Value = "'  +2.74%  "
$ws.Range("D42").Style = $s42D
$ws.Range("E42").Style = $s42E

# Row 43
$s43D = $ws.Range("D43").Style
$s43E = $ws.Range("E43").Style
$ws.Range("D43").Value = "'1.931.18"
$ws.Range("E43").Value = "'  +0.16%  "
$ws.Range("D43").Style = $s43D
$ws.Range("E43").Style = $s43E

# Row 44
$s44D = $ws.Range("D44").Style
$s44E = $ws.Range("E44").Style
$ws.Range("D44").Value = "'0.0279"
$ws.Range("E44").Value = "'  -0.15%  "
$ws.Range("D44").Style = $s44D
$ws.Range("E44").Style = $s44E

# Row 45
$s45E = $ws.Range("E45").Style
$ws.Range("E45").Value = "'  +1.94%  "
$ws.Range("E45").Style = $s45E

# Row 46
$s46E = $ws.Range("E46").Style
$ws.Range("E46").Value = "'  -1.80%  "
$ws.Range("E46").Style = $s46E

# Row 47
$s47D = $ws.Range("D47").Style
$s47E = $ws.Range("E47").Style
$ws.Range("D47").Value = "'9.15"
$ws.Range("E47").Value = "'  -9.55%  "
$ws.Range("D47").Style = $s47D
$ws.Range("E47").Style = $s47E

# Row 48
$s48D = $ws.Range("D48").Style
$s48E = $ws.Range("E48").Style
$ws.Range("D48").Value = "'2.591.21"
$ws.Range("E48").Value = "'  +1.75%  "
$ws.Range("D48").Style = $s48D
$ws.Range("E48").Style = $s48E

# Row 49
$s49E = $ws.Range("E49").Style
$ws.Range("E49").Value = "'  +1.70%  "
$ws.Range("E49").Style = $s49E

# Row 50
$s50D = $ws.Range("D50").Style
$s50E = $ws.Range("E50").Style
$ws.Range("D50").Value = "'71.95"
$ws.Range("E50").Value = "'  -0.50%  "
$ws.Range("D50").Style = $s50D
$ws.Range("E50").Style = $s50E

# Row 51
$s51E = $ws.Range("E51").Style
$ws.Range("E51").Value = "'  +0.70%  "
$ws.Range("E51").Style = $s51E
